# Weekly fruit/vegetable price update: insert 3 new rows of data
# (date 44943, "Región del Maule" origin) at the top of the
# Sandia / Vega Modelo de Temuco block, pushing the existing rows
# 544:649 down to 547:652.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows before row 544, shifting rows 544:649 -> 547:652
$ws.Rows("544:546").Insert()

# --- New row 544: Extra quality ---
$ws.Cells.Item(544,1).Value = 10
$ws.Cells.Item(544,2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(544,3).Value = "La Araucanía"
$ws.Cells.Item(544,4).Value = 44943
$ws.Cells.Item(544,5).Value = 9
$ws.Cells.Item(544,6).Value = 100112028
$ws.Cells.Item(544,7).Value = "Sandia"
$ws.Cells.Item(544,8).Value = "Sin especificar"
$ws.Cells.Item(544,9).Value = "Extra"
$ws.Cells.Item(544,10).Value = 2500
$ws.Cells.Item(544,11).Value = 3500
$ws.Cells.Item(544,12).Value = 3500
$ws.Cells.Item(544,13).Value = 3500
$ws.Cells.Item(544,14).Value = "$/unidad"
$ws.Cells.Item(544,15).Value = "Región del Maule"
$ws.Cells.Item(544,16).Value = 3500
$ws.Cells.Item(544,17).Value = 1
$ws.Cells.Item(544,18).Value = "Hortaliza"

# --- New row 545: Primera quality ---
$ws.Cells.Item(545,1).Value = 10
$ws.Cells.Item(545,2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(545,3).Value = "La Araucanía"
$ws.Cells.Item(545,4).Value = 44943
$ws.Cells.Item(545,5).Value = 9
$ws.Cells.Item(545,6).Value = 100112028
$ws.Cells.Item(545,7).Value = "Sandia"
$ws.Cells.Item(545,8).Value = "Sin especificar"
$ws.Cells.Item(545,9).Value = "Primera"
$ws.Cells.Item(545,10).Value = 3100
$ws.Cells.Item(545,11).Value = 3000
$ws.Cells.Item(545,12).Value = 3000
$ws.Cells.Item(545,13).Value = 3000
$ws.Cells.Item(545,14).Value = "$/unidad"
$ws.Cells.Item(545,15).Value = "Región del Maule"
$ws.Cells.Item(545,16).Value = 3000
$ws.Cells.Item(545,17).Value = 1
$ws.Cells.Item(545,18).Value = "Hortaliza"

# --- New row 546: Segunda quality ---
$ws.Cells.Item(546,1).Value = 10
$ws.Cells.Item(546,2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(546,3).Value = "La Araucanía"
$ws.Cells.Item(546,4).Value = 44943
$ws.Cells.Item(546,5).Value = 9
$ws.Cells.Item(546,6).Value = 100112028
$ws.Cells.Item(546,7).Value = "Sandia"
$ws.Cells.Item(546,8).Value = "Sin especificar"
$ws.Cells.Item(546,9).Value = "Segunda"
$ws.Cells.Item(546,10).Value = 1800
$ws.Cells.Item(546,11).Value = 2500
$ws.Cells.Item(546,12).Value = 2500
$ws.Cells.Item(546,13).Value = 2500
$ws.Cells.Item(546,14).Value = "$/unidad"
$ws.Cells.Item(546,15).Value = "Región del Maule"
$ws.Cells.Item(546,16).Value = 2500
$ws.Cells.Item(546,17).Value = 1
$ws.Cells.Item(546,18).Value = "Hortaliza"

# Apply the same date style (yyyy-mm-dd, style index used by column D)
# as the surrounding cells to the 3 new D-column cells.
$ws.Range("D544:D546").NumberFormat = $ws.Range("D547").NumberFormat
